# Slide 11 ("Today's in-class exercise part B"), Content Placeholder shape (2nd
# shape on the slide). The bullet "Create a new variable with the total GPP
# over time." gets its final sentence expanded with a parenthetical remark
# that contains an ordinal ("1st") with the "st" suffix superscripted.

$p   = $ppt.ActivePresentation
$sl  = $p.Slides.Item(11)
$shp = $sl.Shapes.Item(2)
$tr  = $shp.TextFrame.TextRange

# The bullet is the 2nd paragraph in this placeholder's text body.
$para = $tr.Paragraphs(2, 1)

# Locate "GPP over time." inside that paragraph and replace it with the
# expanded wording (keeping the leading "Create a new variable with the
# total " run/formatting untouched).
$oldTail = "GPP over time."
$start = $para.Text.IndexOf($oldTail) + 1
$len = $oldTail.Length
$tail = $para.Characters($start, $len)

$dash = [char]0x2013
$tail.Text = "GPP over time (not a spatial average $dash just considering the "

# Append the ordinal "1st" and the closing words as their own runs so "st"
# can be superscripted independently of the surrounding text.
[void]$para.InsertAfter("1")
[void]$para.InsertAfter("st")
[void]$para.InsertAfter(" dimension).")

# Superscript just the "st" suffix of the ordinal.
$ordLen = "1st".Length
$stStart = $para.Text.IndexOf("1st") + 1 + 1
$stRange = $para.Characters($stStart, 2)
$stRange.Font.Superscript = $true

$para.Text
